$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.115850925445557
$ws.Range("B1").Value = 3.587602615356445
$ws.Range("C1").Value = 2.833600282669067
$ws.Range("D1").Value = 2.311908721923828
$ws.Range("E1").Value = 1.51703405380249
